# ---------------------------------------------------------------------------
# Edit summary (per the supplied OOXML diff):
#   1. Slide 5's table (graphicFrame) switches its table style from
#      {218C6C99-2EBB-4DC7-9B42-A2E00BC821BD} to
#      {0006199C-1FBE-4705-80B9-054716B92591}.
#   2. The deck's theme ("Integral" / Red Violet colour scheme) is replaced
#      by the stock "Office Theme" colour scheme (font scheme & format
#      scheme are identical between the two themes, so only the 12 theme
#      colours actually change).
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style swap on slide 5 (the graphicFrame holding the table is the
#    2nd shape on that slide).
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{0006199C-1FBE-4705-80B9-054716B92591}")

# ---------------------------------------------------------------------------
# 2) Re-theme the deck: push the stock "Office Theme" colours into the
#    presentation's theme (the 12 scheme slots map 1:1 to
#    dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink).
# ---------------------------------------------------------------------------
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $themeColors.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
